$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look numeric so they remain text
# (matches the source data which stores these as inline strings).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply the updated values
$ws.Range('D2').Value = '26.109.20'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.666.64'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '209.15'
$ws.Range('E5').Value = '  -3.96%  '
$ws.Range('D6').Value = '0.5247'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').Value = '0.2610'
$ws.Range('E8').Value = '  -4.14%  '
$ws.Range('D9').Value = '0.06290'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').Value = '21.05'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('D11').Value = '0.07525'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').Value = '1.675.26'
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('D13').Value = '4.423'
$ws.Range('E13').Value = '  -2.20%  '
$ws.Range('D14').Value = '0.5482'
$ws.Range('E14').Value = '  -5.47%  '
$ws.Range('D15').Value = '66.24'
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('D16').Value = '0.000007941'
$ws.Range('E16').Value = '  -5.16%  '
$ws.Range('D17').Value = '26.128.65'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').Value = '1.002'
$ws.Range('D19').Value = '4.695'
$ws.Range('E19').Value = '  -4.33%  '
$ws.Range('D20').Value = '186.13'
$ws.Range('E20').Value = '  -3.98%  '
$ws.Range('D21').Value = '10.23'
$ws.Range('E21').Value = '  -5.86%  '
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').Value = '1.004'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '149.46'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '0.1240'
$ws.Range('E25').Value = '  -3.59%  '
$ws.Range('D26').Value = '7.450'
$ws.Range('E26').Value = '  -5.39%  '
$ws.Range('D27').Value = '15.86'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '0.06351'
$ws.Range('E28').Value = '  +3.78%  '
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('E30').Value = '  -4.04%  '
$ws.Range('E31').Value = '  -3.16%  '
$ws.Range('D32').Value = '3.405'
$ws.Range('E32').Value = '  -4.98%  '
$ws.Range('D33').Value = '1.631'
$ws.Range('E33').Value = '  -3.37%  '
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('D35').Value = '2.406'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = '0.5983'
$ws.Range('E36').Value = '  -3.38%  '
$ws.Range('D37').Value = '2.742'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = '6.106'
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.105.62'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('D41').Value = '0.8694'
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('D43').Value = '99.84'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('D44').Value = '1.817.88'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = '0.00000000107'
$ws.Range('E45').Value = '  -4.01%  '
$ws.Range('D46').Value = '55.16'
$ws.Range('E46').Value = '  -4.59%  '
$ws.Range('D47').Value = '0.9978'
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').Value = '7.992'
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('D49').Value = '0.05226'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').Value = '0.4244'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').Value = '5.912'
$ws.Range('E51').Value = '  -2.40%  '
